$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44699
$ws.Range('L2').Value2 = 'Primera'
$ws.Range("M2").Value2 = 100
$ws.Range("N2").Value2 = 20000
$ws.Range("O2").Value2 = 22000
$ws.Range("P2").Value2 = 21000
$ws.Range('Q2').Value2 = '$/caja 18 kilos'
$ws.Range("S2").Value2 = 1167

# Row 3
$ws.Range("D3").Value2 = 44699
$ws.Range('L3').Value2 = 'Segunda'
$ws.Range("M3").Value2 = 50
$ws.Range("N3").Value2 = 18000
$ws.Range("O3").Value2 = 18000
$ws.Range("P3").Value2 = 18000
$ws.Range('Q3').Value2 = '$/caja 18 kilos'
$ws.Range("S3").Value2 = 1000

# Row 4
$ws.Range("D4").Value2 = 45084
$ws.Range('L4').Value2 = 'Primera'
$ws.Range("M4").Value2 = 100
$ws.Range("N4").Value2 = 20000
$ws.Range("O4").Value2 = 21000
$ws.Range("P4").Value2 = 20500
$ws.Range('Q4').Value2 = '$/caja 18 kilos granel'
$ws.Range("S4").Value2 = 1139

# Row 5
$ws.Range("D5").Value2 = 45002
$ws.Range('L5').Value2 = 'Primera'
$ws.Range("M5").Value2 = 100
$ws.Range("N5").Value2 = 12000
$ws.Range("O5").Value2 = 13000
$ws.Range("P5").Value2 = 12500
$ws.Range('Q5').Value2 = '$/caja 18 kilos'
$ws.Range("S5").Value2 = 694

# Row 6
$ws.Range("D6").Value2 = 45014
$ws.Range('L6').Value2 = 'Primera'
$ws.Range("M6").Value2 = 50
$ws.Range("N6").Value2 = 13000
$ws.Range("O6").Value2 = 14000
$ws.Range("P6").Value2 = 13600
$ws.Range('Q6').Value2 = '$/caja 18 kilos'
$ws.Range("S6").Value2 = 756

# Row 7
$ws.Range("D7").Value2 = 45014
$ws.Range('L7').Value2 = 'Segunda'
$ws.Range("M7").Value2 = 20
$ws.Range("N7").Value2 = 10000
$ws.Range("O7").Value2 = 10000
$ws.Range("P7").Value2 = 10000
$ws.Range('Q7').Value2 = '$/caja 18 kilos'
$ws.Range("S7").Value2 = 556

# Row 8
$ws.Range("D8").Value2 = 45044
$ws.Range('L8').Value2 = 'Primera'
$ws.Range("M8").Value2 = 100
$ws.Range("N8").Value2 = 17000
$ws.Range("O8").Value2 = 18000
$ws.Range("P8").Value2 = 17500
$ws.Range('Q8').Value2 = '$/caja 18 kilos'
$ws.Range("S8").Value2 = 972

# Row 9
$ws.Range("D9").Value2 = 45030
$ws.Range('L9').Value2 = 'Primera'
$ws.Range("M9").Value2 = 100
$ws.Range("N9").Value2 = 15000
$ws.Range("O9").Value2 = 16000
$ws.Range("P9").Value2 = 15500
$ws.Range('Q9').Value2 = '$/caja 18 kilos granel'
$ws.Range("S9").Value2 = 861

# Row 10
$ws.Range("D10").Value2 = 44819
$ws.Range('L10').Value2 = 'Primera'
$ws.Range("M10").Value2 = 100
$ws.Range("N10").Value2 = 25000
$ws.Range("O10").Value2 = 26000
$ws.Range("P10").Value2 = 25500
$ws.Range('Q10').Value2 = '$/caja 18 kilos granel'
$ws.Range("S10").Value2 = 1417

# Row 11
$ws.Range("D11").Value2 = 44280
$ws.Range('L11').Value2 = 'Primera'
$ws.Range("M11").Value2 = 100
$ws.Range("N11").Value2 = 14000
$ws.Range("O11").Value2 = 15000
$ws.Range("P11").Value2 = 14500
$ws.Range('Q11').Value2 = '$/caja 18 kilos'
$ws.Range("S11").Value2 = 806

# Row 12
$ws.Range("D12").Value2 = 44280
$ws.Range('L12').Value2 = 'Segunda'
$ws.Range("M12").Value2 = 50
$ws.Range("N12").Value2 = 12000
$ws.Range("O12").Value2 = 12000
$ws.Range("P12").Value2 = 12000
$ws.Range('Q12').Value2 = '$/caja 18 kilos'
$ws.Range("S12").Value2 = 667

# Row 13
$ws.Range("D13").Value2 = 44516
$ws.Range('L13').Value2 = 'Primera'
$ws.Range("M13").Value2 = 100
$ws.Range("N13").Value2 = 33000
$ws.Range("O13").Value2 = 34000
$ws.Range("P13").Value2 = 33500
$ws.Range('Q13').Value2 = '$/caja 18 kilos'
$ws.Range("S13").Value2 = 1861

# Row 14
$ws.Range("D14").Value2 = 44316
$ws.Range('L14').Value2 = 'Primera'
$ws.Range("M14").Value2 = 50
$ws.Range("N14").Value2 = 20000
$ws.Range("O14").Value2 = 20000
$ws.Range("P14").Value2 = 20000
$ws.Range('Q14').Value2 = '$/caja 18 kilos'
$ws.Range("S14").Value2 = 1111

